$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: A4 loses its special "keyword" font, becomes plain default text ---
$ws.Range("A4").Font.Name = "Nachlieli CLM"

# --- Row 5: A5 and B5 become plain default text too ---
$ws.Range("A5").Font.Name = "Nachlieli CLM"
$ws.Range("B5").Font.Name = "Nachlieli CLM"

# --- Row 6: keyword text gets padded with spaces, wraps, row taller ---
$ws.Range("A6").Value = " free games helix jump "
$ws.Range("A6").Font.Name = "Nachlieli CLM"
$ws.Range("A6").WrapText = $true
$ws.Range("B6").Font.Name = "Nachlieli CLM"
$ws.Rows.Item(6).RowHeight = 24

# C6 / D6 become real hyperlinks (were empty placeholder cells before),
# styled with the workbook's existing blue hyperlink-like font.
$ws.Hyperlinks.Add($ws.Range("C6"), "https://play.google.com/store/apps/details?id=com.singleton.helix", "", "", "Explore") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "https://play.google.com/store/apps/details?id=com.singleton.helix", "", "", "Show (4)") | Out-Null
$ws.Range("C6:D6").Font.Name = "Nachlieli CLM"
$ws.Range("C6:D6").Font.Color = 16711680
$ws.Range("C6:D6").WrapText = $true

# --- Row 7: replaced with a duplicate of row 5's keyword pair ---
$ws.Range("A7").Value = "helix jump"
$ws.Range("B7").Value = "com.singleton.helix"
$ws.Range("A7").Font.Name = "Nachlieli CLM"
$ws.Range("B7").Font.Name = "Nachlieli CLM"
$ws.Rows.Item(7).RowHeight = 12.8

# --- Column A narrower ---
$ws.Columns.Item(1).ColumnWidth = 20.57

# --- Selection moves ---
$ws.Range("B10").Select()
